# Day 31 content updated

$p = $ppt.ActivePresentation

# --- Slide 1: title "DAY 30" -> "DAY 32" ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "DAY 32"

# --- Slide 2: agenda bullet list updates ---
$s2 = $p.Slides.Item(2)
$shape = $s2.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph "PV and pvc demo" -> split the leading "PV and " run into
# "pv" + " and " (the trailing "pvc" / " demo" runs stay untouched).
# There are two "PV and " occurrences in this text box, so skip past the
# first one ("PV and PVC Data dependencies...") to reach the right one.
$firstPvAnd = $tr.Find("PV and ", 0, $true)
$pvAnd = $tr.Find("PV and ", $firstPvAnd.Start + $firstPvAnd.Length, $true)
$pvOnly = $tr.Characters($pvAnd.Start, 2)
$pvOnly.Text = "pv"

# Paragraph "statefulset understanding and use cases" -> "statefulSet ..."
$sf1 = $tr.Find("statefulset", 0, $true)
$sf1.Text = "statefulSet"

# Paragraph "Statefulset Demo" -> "StatefulSet Demo"
$sf2 = $tr.Find("Statefulset", 0, $true)
$sf2.Text = "StatefulSet"
